$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.747.36"
$ws.Range("E2").Value = "  -3.47%  "
$ws.Range("D3").Value = "2.487.83"
$ws.Range("E3").Value = "  -6.05%  "
$ws.Range("E4").Value = "  +0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "556.73"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -4.33%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "148.67"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -4.93%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -3.25%  "
$ws.Range("D9").Value = "2.486.41"
$ws.Range("E9").Value = "  -5.99%  "
$ws.Range("E10").Value = "  -7.91%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.50"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -5.22%  "
$ws.Range("E12").Value = "  -1.25%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.360"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -6.16%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "26.54"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -7.31%  "
$ws.Range("D15").Value = "2.936.31"
$ws.Range("E15").Value = "  -5.92%  "
$ws.Range("E16").Value = "  -8.62%  "
$ws.Range("D17").Value = "61.622.08"
$ws.Range("E17").Value = "  -3.51%  "
$ws.Range("D18").Value = "2.484.30"
$ws.Range("E18").Value = "  -5.77%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.29"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.16"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -7.52%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.24"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -6.79%  "
$ws.Range("E22").Value = "  -6.43%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +3.08%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "64.41"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -5.30%  "
$ws.Range("E26").Value = "  -9.12%  "
$ws.Range("D27").Value = "2.601.74"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "555.39"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -7.64%  "
$ws.Range("E29").Value = "  -5.48%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.93"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.09%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.40"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -9.59%  "
$ws.Range("E33").Value = "  -5.36%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.93"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -7.01%  "
$ws.Range("E35").Value = "  -7.69%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "6.01"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -9.68%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "4.97"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -9.40%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  -4.58%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "18.69"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -5.47%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "146.96"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("E42").Value = "  -7.06%  "
$ws.Range("E43").Value = "  +0.06%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "40.60"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.18%  "
$ws.Range("E45").Value = "  -4.52%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "148.74"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -9.18%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "22.21"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -8.15%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "3.68"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -5.91%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0547"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -7.45%  "
$ws.Range("E50").Value = "  -5.52%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0949"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -5.14%  "
